# Update the "Marketplace Report" sheet with the latest generated
# notification data (add a function for notification).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Top summary block -------------------------------------------------
# "Generated on:" timestamp
$ws.Range("B1").Value = "March 22, 2025 at 03:34:21 PM"

# "Total Sales (P):" summary (kept as formatted text, like the source,
# not a plain number, so force text entry then restore the default format)
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "3,423"
$ws.Range("B3").ClearFormats()

# "Total Auction Sessions:" summary
$ws.Range("B4").Value = 1

# --- Monthly breakdown table --------------------------------------------
# February row (row 8): Total Users
$ws.Range("C8").Value = 1

# March row (row 9): Total Sales / Total Users / Total Auction Sessions
$ws.Range("B9").Value = 3423
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 1

# TOTAL row (row 20): Total Sales / Total Auction Sessions
$ws.Range("B20").Value = 3423
$ws.Range("D20").Value = 1
